$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 691.6667
$ws.Range("I6").Value = 37.166668
$ws.Range("J6").Value = 2000.6666
$ws.Range("K6").Value = 111.500004
$ws.Range("L6").Value = 6001.9998
$ws.Range("M6").Value = 0.4999959999999959
$ws.Range("N6").Value = -6225.9998
$ws.Range("H31").Value = 4275.75
$ws.Range("J31").Value = 8500
$ws.Range("L31").Value = 25500
$ws.Range("N31").Value = -25960
$ws.Range("H33").Value = 88.125
$ws.Range("I33").Value = 98.416664
$ws.Range("J33").Value = 57.25
$ws.Range("K33").Value = 98.416664
$ws.Range("L33").Value = 57.25
$ws.Range("M33").Value = 130.583336
$ws.Range("N33").Value = -515.25
$ws.Range("H99").Value = 1875.3
$ws.Range("I99").Value = 570.5
$ws.Range("K99").Value = 1711.5
$ws.Range("M99").Value = -213.5
$ws.Range("H132").Value = 13440
$ws.Range("I132").Value = 13375
$ws.Range("K132").Value = 40125
$ws.Range("M132").Value = -37595
$ws.Range("H141").Value = 2610.8
$ws.Range("I141").Value = 2321.3333
$ws.Range("J141").Value = 3045
$ws.Range("K141").Value = 6963.999899999999
$ws.Range("L141").Value = 9135
$ws.Range("M141").Value = -1783.999899999999
$ws.Range("N141").Value = -19495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 1850
$ws.Range("I10").Value = 1850
$ws.Range("K10").Value = 1850
$ws.Range("M10").Value = -1680
$ws.Range("H32").Value = 4403.032
$ws.Range("I32").Value = 4282.6665
$ws.Range("K32").Value = 4282.6665
$ws.Range("M32").Value = -3995.6665
$ws.Range("H61").Value = 5796.5386
$ws.Range("I61").Value = 5388.6665
$ws.Range("K61").Value = 5388.6665
$ws.Range("M61").Value = -5176.6665
$ws.Range("H63").Value = 2379.3
$ws.Range("J63").Value = 5000
$ws.Range("L63").Value = 5000
$ws.Range("N63").Value = -6372
$ws.Range("H66").Value = 2379.3
$ws.Range("J66").Value = 5000
$ws.Range("L66").Value = 25000
$ws.Range("N66").Value = -31864
$ws.Range("H136").Value = 5796.5386
$ws.Range("I136").Value = 5388.6665
$ws.Range("K136").Value = 16165.9995
$ws.Range("M136").Value = -13615.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 14073.625
$ws.Range("I82").Value = 14073.625
$ws.Range("K82").Value = 14073.625
$ws.Range("M82").Value = -13690.625
$ws.Range("H85").Value = 14073.625
$ws.Range("I85").Value = 14073.625
$ws.Range("K85").Value = 14073.625
$ws.Range("M85").Value = -12747.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 373.8889
$ws.Range("I35").Value = 374.375
$ws.Range("K35").Value = 374.375
$ws.Range("M35").Value = -80.375
$ws.Range("H134").Value = 3524.3333
$ws.Range("I134").Value = 2388.4285
$ws.Range("J134").Value = 7500
$ws.Range("K134").Value = 7165.2855
$ws.Range("L134").Value = 22500
$ws.Range("M134").Value = -4630.2855
$ws.Range("N134").Value = -27570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 200447.2
$ws.Range("I11").Value = 334015.34
$ws.Range("J11").Value = 95
$ws.Range("K11").Value = 1002046.02
$ws.Range("L11").Value = 285
$ws.Range("M11").Value = -1001906.02
$ws.Range("N11").Value = -565
$ws.Range("H12").Value = 111.4
$ws.Range("I12").Value = 133.5
$ws.Range("J12").Value = 105.875
$ws.Range("K12").Value = 400.5
$ws.Range("L12").Value = 317.625
$ws.Range("M12").Value = -227.5
$ws.Range("N12").Value = -663.625
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("L22").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("L27").Value = 0
$ws.Range("H109").Value = 102763.6
$ws.Range("I109").Value = 144779.72
$ws.Range("J109").Value = 4726
$ws.Range("K109").Value = 434339.16
$ws.Range("L109").Value = 14178
$ws.Range("M109").Value = -433299.16
$ws.Range("N109").Value = -16258
$ws.Range("H137").Value = 2449.5
$ws.Range("I137").Value = 1000
$ws.Range("K137").Value = 3000
$ws.Range("M137").Value = 2100

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 27399.6
$ws.Range("I26").Value = 25000
$ws.Range("J26").Value = 27999.5
$ws.Range("K26").Value = 25000
$ws.Range("L26").Value = 27999.5
$ws.Range("M26").Value = -24720
$ws.Range("N26").Value = -28559.5
$ws.Range("H46").Value = 9454.5
$ws.Range("J46").Value = 9614.666999999999
$ws.Range("L46").Value = 9614.666999999999
$ws.Range("N46").Value = -9926.666999999999
$ws.Range("H50").Value = 27399.6
$ws.Range("I50").Value = 25000
$ws.Range("J50").Value = 27999.5
$ws.Range("K50").Value = 25000
$ws.Range("L50").Value = 27999.5
$ws.Range("M50").Value = -24502
$ws.Range("N50").Value = -28995.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 24499.334
$ws.Range("I56").Value = 23999.5
$ws.Range("K56").Value = 23999.5
$ws.Range("M56").Value = -23308.5
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("L76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("L79").Value = 0
$ws.Range("H100").Value = 6946.8276
$ws.Range("I100").Value = 4868.1
$ws.Range("J100").Value = 8040.8945
$ws.Range("K100").Value = 4868.1
$ws.Range("L100").Value = 8040.8945
$ws.Range("M100").Value = -4327.1
$ws.Range("N100").Value = -9122.8945
$ws.Range("H123").Value = 89998
$ws.Range("J123").Value = 89998
$ws.Range("L123").Value = 89998
$ws.Range("N123").Value = -99798

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 22249.75
$ws.Range("I58").Value = 14999.5
$ws.Range("J58").Value = 29500
$ws.Range("K58").Value = 14999.5
$ws.Range("L58").Value = 29500
$ws.Range("M58").Value = -14691.5
$ws.Range("N58").Value = -30116
$ws.Range("H61").Value = 6016.6665
$ws.Range("I61").Value = 6016.6665
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6016.6665
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("M61").Value = -5724.6665
